$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 4043.2222
$ws.Range("I131").Value = 4541.2856
$ws.Range("J131").Value = 2300
$ws.Range("K131").Value = 13623.8568
$ws.Range("L131").Value = 6900
$ws.Range("M131").Value = -8583.856800000001
$ws.Range("N131").Value = -16980
$ws.Range("H137").Value = 5560608.5
$ws.Range("I137").Value = 12508175
$ws.Range("K137").Value = 37524525
$ws.Range("M137").Value = -37521975

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3526.7334
$ws.Range("I61").Value = 1492.1
$ws.Range("J61").Value = 7596
$ws.Range("K61").Value = 1492.1
$ws.Range("L61").Value = 7596
$ws.Range("M61").Value = -1280.1
$ws.Range("N61").Value = -8020
$ws.Range("H74").Value = 1364.25
$ws.Range("I74").Value = 1159.8
$ws.Range("J74").Value = 1705
$ws.Range("K74").Value = 1159.8
$ws.Range("L74").Value = 1705
$ws.Range("M74").Value = -285.8
$ws.Range("N74").Value = -3453
$ws.Range("H77").Value = 1364.25
$ws.Range("I77").Value = 1159.8
$ws.Range("J77").Value = 1705
$ws.Range("K77").Value = 5799
$ws.Range("L77").Value = 8525
$ws.Range("M77").Value = -1431
$ws.Range("N77").Value = -17261
$ws.Range("H80").Value = 28795
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 28795
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 28795
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -30791
$ws.Range("H83").Value = 28795
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 28795
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 86385
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -96369
$ws.Range("H132").Value = 37042270
$ws.Range("I132").Value = 43483320
$ws.Range("K132").Value = 130449960
$ws.Range("M132").Value = -130447430
$ws.Range("H136").Value = 3526.7334
$ws.Range("I136").Value = 1492.1
$ws.Range("J136").Value = 7596
$ws.Range("K136").Value = 4476.299999999999
$ws.Range("L136").Value = 22788
$ws.Range("M136").Value = -1926.299999999999
$ws.Range("N136").Value = -27888

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 31669
$ws.Range("I15").Value = 12500
$ws.Range("J15").Value = 70007
$ws.Range("K15").Value = 12500
$ws.Range("L15").Value = 70007
$ws.Range("M15").Value = -12273
$ws.Range("N15").Value = -70461
$ws.Range("H80").Value = 1074.4546
$ws.Range("I80").Value = 1193.2222
$ws.Range("K80").Value = 1193.2222
$ws.Range("M80").Value = -195.2221999999999
$ws.Range("H83").Value = 1074.4546
$ws.Range("I83").Value = 1193.2222
$ws.Range("K83").Value = 5966.111
$ws.Range("M83").Value = -974.1109999999999
$ws.Range("H134").Value = 2236.742
$ws.Range("I134").Value = 1901.6316
$ws.Range("J134").Value = 2767.3333
$ws.Range("K134").Value = 5704.8948
$ws.Range("L134").Value = 8301.999899999999
$ws.Range("M134").Value = -3169.8948
$ws.Range("N134").Value = -13371.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3228177.5
$ws.Range("H34").Value = 3228177.5
$ws.Range("H50").Value = 16000
$ws.Range("J50").Value = 16000
$ws.Range("L50").Value = 16000
$ws.Range("N50").Value = -17250
$ws.Range("H51").Value = 1000000000
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H60").Value = 11659.538
$ws.Range("I60").Value = 10000
$ws.Range("J60").Value = 11797.833
$ws.Range("K60").Value = 10000
$ws.Range("L60").Value = 11797.833
$ws.Range("M60").Value = -9489
$ws.Range("N60").Value = -12819.833
$ws.Range("H61").Value = 1000000000
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H86").Value = 6283.75
$ws.Range("I86").Value = 5501.25
$ws.Range("J86").Value = 6675
$ws.Range("K86").Value = 5501.25
$ws.Range("L86").Value = 6675
$ws.Range("M86").Value = -4378.25
$ws.Range("N86").Value = -8921
$ws.Range("H89").Value = 6283.75
$ws.Range("I89").Value = 5501.25
$ws.Range("J89").Value = 6675
$ws.Range("K89").Value = 27506.25
$ws.Range("L89").Value = 33375
$ws.Range("M89").Value = -21890.25
$ws.Range("N89").Value = -44607
$ws.Range("H122").Value = 2084
$ws.Range("I122").Value = 1779.25
$ws.Range("J122").Value = 2541.125
$ws.Range("K122").Value = 5337.75
$ws.Range("L122").Value = 7623.375
$ws.Range("M122").Value = -2887.75
$ws.Range("N122").Value = -12523.375
$ws.Range("H134").Value = 2378.2
$ws.Range("I134").Value = 1656.7059
$ws.Range("J134").Value = 6466.6665
$ws.Range("K134").Value = 4970.1177
$ws.Range("L134").Value = 19399.9995
$ws.Range("M134").Value = -2435.1177
$ws.Range("N134").Value = -24469.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 743.375
$ws.Range("I131").Value = 439.16666
$ws.Range("J131").Value = 1290.95
$ws.Range("K131").Value = 1317.49998
$ws.Range("L131").Value = 3872.85
$ws.Range("M131").Value = 3722.50002
$ws.Range("N131").Value = -13952.85

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 523689.75
$ws.Range("I126").Value = 1742.7142
$ws.Range("J126").Value = 738609.1
$ws.Range("K126").Value = 5228.142599999999
$ws.Range("L126").Value = 2215827.3
$ws.Range("M126").Value = -2758.142599999999
$ws.Range("N126").Value = -2220767.3
$ws.Range("H132").Value = 1941.4482
$ws.Range("I132").Value = 1532.08
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 4596.24
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -2066.24
$ws.Range("N132").Value = -18560

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1991.4073
$ws.Range("I132").Value = 1587.1111
$ws.Range("J132").Value = 2800
$ws.Range("K132").Value = 4761.3333
$ws.Range("L132").Value = 8400
$ws.Range("M132").Value = -2231.3333
$ws.Range("N132").Value = -13460

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 620
$ws.Range("I81").Value = 584.61536
$ws.Range("K81").Value = 1169.23072
$ws.Range("M81").Value = -108.23072
$ws.Range("H82").Value = 30014.285
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 30014.285
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 30014.285
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -30780.285
$ws.Range("H84").Value = 620
$ws.Range("I84").Value = 584.61536
$ws.Range("K84").Value = 5846.1536
$ws.Range("M84").Value = -542.1535999999996
$ws.Range("H85").Value = 30014.285
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 30014.285
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 30014.285
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -32666.285
$ws.Range("H132").Value = 317951.53
$ws.Range("I132").Value = 386601.88
$ws.Range("J132").Value = 20466.666
$ws.Range("K132").Value = 1159805.64
$ws.Range("L132").Value = 61399.99800000001
$ws.Range("M132").Value = -1157275.64
$ws.Range("N132").Value = -66459.99800000001
